$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right after row 436 (before the old row 437),
# pushing the existing rows 437:547 down to 438:548.
$ws.Rows("437:437").Insert()

# The row that used to be 437 is now at 438; duplicate it into the freshly
# inserted (blank) row 437 so every column but the date matches the diff.
$ws.Range("A438:R438").Copy()
$ws.Range("A437:R437").PasteSpecial()

# The new record's date is the only field that differs from the row it was
# copied from.
$ws.Range("D437").Value = [DateTime]"2023-07-28"
